$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Veronia Rafat, Administrator, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G3").Value = 'Dr. Veronia Rafat, Dr. Hend Mahmoud, Administrator, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Majorelle Magdy'
$ws.Range("G4").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Majorelle Magdy'
$ws.Range("G5").Value = 'Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda'
$ws.Range("G6").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Mohammad El-Tanany, Dr. Majorelle Magdy, Dr. Alshimaa Atef, Dr. Manar Montaser'
$ws.Range("G7").Value = 'Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G8").Value = 'Dr. Nada Mohammad, Dr. Abeer Ragab'
$ws.Range("G11").Value = 'Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range("G12").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Madeha Saeed'
$ws.Range("G13").Value = 'Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh'
$ws.Range("G24").Value = 'Dr. Sarah Mahdy, Dr. Youstina Gamil'
$ws.Range("G27").Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Range("G28").Value = 'Dr. Aya Emad, Dr. Maryam Ashraf'
$ws.Range("G30").Value = 'Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Shorok Mohammad'
